$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume refresh. For the "Price" column, some new values look
# like plain decimals (e.g. "1.001", "0.3817") which Excel would otherwise
# auto-convert to numbers (dropping significant trailing zeros etc.), so each
# such cell is forced to Text format before the literal string is written.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.993.24"
$ws.Range("E2").Value = "  -1.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.882.70"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.34"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4999"
$ws.Range("E7").Value = "  -3.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3817"
$ws.Range("E8").Value = "  -3.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09171"
$ws.Range("E9").Value = "  -5.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.120"
$ws.Range("E10").Value = "  -2.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.65"
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.316"
$ws.Range("E12").Value = "  -3.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.72"
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.883.88"
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.267"
$ws.Range("E15").Value = "  -3.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001104"
$ws.Range("E17").Value = "  -2.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.35"
$ws.Range("E18").Value = "  -3.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06642"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.93"
$ws.Range("E20").Value = "  -1.56%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.172"
$ws.Range("E22").Value = "  -2.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.057.56"
$ws.Range("E23").Value = "  -1.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.39"
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.298"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.100.20"
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.542"
$ws.Range("E27").Value = "  -5.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.78"
$ws.Range("E28").Value = "  -2.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "157.18"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.53"
$ws.Range("E30").Value = "  -1.76%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1053"
$ws.Range("E31").Value = "  -2.34%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.062"
$ws.Range("E32").Value = "  -4.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.582"
$ws.Range("E33").Value = "  -3.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.592"
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.332"
$ws.Range("E35").Value = "  -7.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06570"
$ws.Range("E36").Value = "  -3.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02409"
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2201"
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.286"
$ws.Range("E39").Value = "  +8.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.207"
$ws.Range("E40").Value = "  -5.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6408"
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.54"
$ws.Range("E42").Value = "  -2.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.925"
$ws.Range("E43").Value = "  -3.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.29"
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6031"
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.285"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.671"
$ws.Range("E48").Value = "  -2.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.988"
$ws.Range("E49").Value = "  -2.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.211"
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.47"
$ws.Range("E51").Value = "  -3.68%  "
